$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Push door" column header to "Force door" (also renames the Table2 column).
$ws.Range("J2").Value = "Force door"

# Update the data: row 5's "Force door" action is now scored as true (1).
$ws.Range("J5").Value = 1

# Add the new "Violent Ending" expected outcome text for row 5.
$ws.Range("R5").Value = "Violent Ending"

# Move the active selection to Q5, matching the author's final cursor position.
$ws.Activate() | Out-Null
$ws.Range("Q5").Select() | Out-Null
